# Update the four report sheets (日间医生/日间科室/预住院医生/预住院科室) with the
# refreshed doctor/department counts - day-surgery totals rose from 38 to 49
# doctors (11 -> 13 departments, 691 -> 817 patients) and pre-admission
# totals rose from 61 to 71 doctors (15 -> 16 departments, 1352 -> 1610
# patients). Each row is rewritten (label + count), sorted descending by
# count, with the trailing "total doctors"/"total patients" summary rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$rows1 = @(
    @('主管医生', '49位医生办理日间手术', 1),
    @('吴华荣', 115, 0),
    @('王会旺', 54, 0),
    @('徐伟坤', 51, 0),
    @('蔡金生', 50, 0),
    @('张梦鑫', 50, 0),
    @('刘少川', 43, 0),
    @('杨克', 38, 0),
    @('张尚普', 35, 0),
    @('白栩搏', 30, 0),
    @('陈子奇', 30, 0),
    @('潘来辉', 29, 0),
    @('贾丁丁', 28, 0),
    @('张连锁', 27, 0),
    @('刘阳', 27, 0),
    @('孙亚东', 25, 0),
    @('刘学沛', 21, 0),
    @('宁胜华', 20, 0),
    @('游小军', 18, 0),
    @('张昕', 17, 0),
    @('李琰', 16, 0),
    @('范会龙', 14, 0),
    @('李继凯', 9, 0),
    @('郝运兵', 8, 0),
    @('王少锋', 7, 0),
    @('白杰', 7, 0),
    @('刘庆辉', 6, 0),
    @('王学攀', 4, 0),
    @('郭二松', 4, 0),
    @('张翼飞', 4, 0),
    @('赵书明', 4, 0),
    @('魏召劝', 3, 0),
    @('李少青', 2, 0),
    @('刘月星', 2, 0),
    @('刘林周', 2, 0),
    @('薛晓乐', 2, 0),
    @('高文华', 2, 0),
    @('朱文博', 1, 0),
    @('张辰阳', 1, 0),
    @('曲巧格', 1, 0),
    @('王振辉', 1, 0),
    @('赵鹏浩', 1, 0),
    @('高勇岗', 1, 0),
    @('程子文', 1, 0),
    @('刘学亮', 1, 0),
    @('段智睿', 1, 0),
    @('刘玉波', 1, 0),
    @('王冬月', 1, 0),
    @('申军国', 1, 0),
    @('程旭', 1, 0),
    @('医生总数', 49, 0),
    @('总计病人', 817, 0)
)
$r = 1
foreach ($row in $rows1) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[2] -eq 1) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 2).Value = [double]$row[1]
    }
    $r++
}

$ws = $wb.Worksheets.Item(2)
$rows2 = @(
    @('负责科室', '13个科室办理日间手术', 1),
    @('骨六科', 365, 0),
    @('骨九科', 204, 0),
    @('外三科', 136, 0),
    @('骨十一科', 39, 0),
    @('骨七科', 35, 0),
    @('骨一科', 10, 0),
    @('骨二科', 8, 0),
    @('骨三科', 7, 0),
    @('骨四科', 4, 0),
    @('骨五科', 3, 0),
    @('外一科', 2, 0),
    @('骨八科', 2, 0),
    @('骨十科', 2, 0),
    @('科室总数', 13, 0),
    @('总计病人', 817, 0)
)
$r = 1
foreach ($row in $rows2) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[2] -eq 1) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 2).Value = [double]$row[1]
    }
    $r++
}

$ws = $wb.Worksheets.Item(3)
$rows3 = @(
    @('主管医生', '71位医生办理预住院', 1),
    @('吴华荣', 122, 0),
    @('曹旭阳', 106, 0),
    @('魏召劝', 98, 0),
    @('张梦鑫', 93, 0),
    @('张尚普', 70, 0),
    @('张翼飞', 69, 0),
    @('贾科锋', 63, 0),
    @('陈子奇', 61, 0),
    @('杨克', 61, 0),
    @('王会旺', 57, 0),
    @('徐伟坤', 54, 0),
    @('蔡金生', 53, 0),
    @('贾丁丁', 53, 0),
    @('孔涛涛', 52, 0),
    @('陈润', 47, 0),
    @('刘阳', 47, 0),
    @('刘少川', 45, 0),
    @('张辰阳', 34, 0),
    @('孙亚东', 31, 0),
    @('潘来辉', 30, 0),
    @('李琰', 30, 0),
    @('白栩搏', 30, 0),
    @('张连锁', 30, 0),
    @('宁胜华', 24, 0),
    @('郎彦飞', 23, 0),
    @('刘学沛', 22, 0),
    @('游小军', 19, 0),
    @('张昕', 18, 0),
    @('范会龙', 18, 0),
    @('郝运兵', 17, 0),
    @('李继凯', 9, 0),
    @('王少锋', 9, 0),
    @('郭二松', 8, 0),
    @('左百军', 8, 0),
    @('刘庆辉', 7, 0),
    @('白杰', 7, 0),
    @('赵玉龙', 7, 0),
    @('王学攀', 7, 0),
    @('刘玉波', 7, 0),
    @('王彦伟', 6, 0),
    @('吕庆列', 5, 0),
    @('刘月星', 5, 0),
    @('赵书明', 5, 0),
    @('李向科', 3, 0),
    @('李栋', 3, 0),
    @('李少青', 3, 0),
    @('曲巧格', 3, 0),
    @('侯心昕', 3, 0),
    @('薛晓乐', 2, 0),
    @('李德磊', 2, 0),
    @('高文华', 2, 0),
    @('段智睿', 2, 0),
    @('刘林周', 2, 0),
    @('朱文博', 1, 0),
    @('范永强', 1, 0),
    @('王洪庆', 1, 0),
    @('刘学亮', 1, 0),
    @('王振辉', 1, 0),
    @('杨良栋', 1, 0),
    @('赵鹏浩', 1, 0),
    @('何举仁', 1, 0),
    @('李强', 1, 0),
    @('陈国江', 1, 0),
    @('高少科', 1, 0),
    @('程子文', 1, 0),
    @('郭志刚', 1, 0),
    @('高勇岗', 1, 0),
    @('孙国栋', 1, 0),
    @('王冬月', 1, 0),
    @('申军国', 1, 0),
    @('程旭', 1, 0),
    @('医生总数', 71, 0),
    @('总计病人', 1610, 0)
)
$r = 1
foreach ($row in $rows3) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[2] -eq 1) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 2).Value = [double]$row[1]
    }
    $r++
}

$ws = $wb.Worksheets.Item(4)
$rows4 = @(
    @('负责科室', '16个科室办理预住院', 1),
    @('骨二科', 469, 0),
    @('骨六科', 401, 0),
    @('骨九科', 372, 0),
    @('外三科', 141, 0),
    @('骨七科', 81, 0),
    @('骨十一科', 43, 0),
    @('骨十二科', 38, 0),
    @('骨一科', 26, 0),
    @('骨四科', 11, 0),
    @('骨三科', 7, 0),
    @('骨十科', 5, 0),
    @('外一科', 5, 0),
    @('骨五科', 4, 0),
    @('骨十五科', 3, 0),
    @('骨八科', 3, 0),
    @('外四科', 1, 0),
    @('科室总数', 16, 0),
    @('总计病人', 1610, 0)
)
$r = 1
foreach ($row in $rows4) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[2] -eq 1) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 2).Value = [double]$row[1]
    }
    $r++
}
